$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dims = $ws.UsedRange
$lastRow = $dims.Rows.Count - 1
if ($lastRow -lt 2) { $lastRow = 173 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
